# Update automatico via Actualizar 03-09-2021 13-05-03
#
# The "D" column holds the timestamp of each availability check. Every run
# of the updater shifts the history down by one slot: the newest batch
# (rows 2-15) receives the freshly captured timestamp, while the values
# that used to occupy rows 2-15 / 16-29 slide down into rows 16-29 / 30-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newest = 44264.54492931246
$mid    = 44264.52353865741
$oldest = 44264.50214331019

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newest
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $mid
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldest
}
